$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header column from "year_group" to "group" (also updates the table column name)
$ws.Range("B1").Value = "group"

# Add new data row values
$ws.Range("A3").Value = "test@gmail.com"
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 2

# Add hyperlink on the new email cell
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:test@gmail.com")

# Copy formatting (Hyperlink cell style) from the row above onto the new row
$ws.Range("A2:C2").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Expand the table to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C3"))

# Update selection to match target state
$ws.Range("F4").Select()
